$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.991.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.745.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.20%  '

$ws.Range("E6").Value = '  +0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5142'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.08%  '

$ws.Range("E8").Value = '  -1.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.740.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.46%  '

$ws.Range("E11").Value = '  +1.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.20'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.16%  '

$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.634'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.10%  '

$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.023.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006822'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.962.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.299'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.08%  '

$ws.Range("E23").Value = '  -1.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.358'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.34%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.504'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.27'
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = '  -1.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.955'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08252'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.673'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04676'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.657'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6241'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.17%  '

$ws.Range("E37").Value = '  +1.08%  '

$ws.Range("E38").Value = '  +0.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.930'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.96%  '

$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.42'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3882'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.77%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7557'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.69%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.021'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.351'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1133'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05234'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.69'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.76%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.633'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3442'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.92%  '
